$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header year column (N3 = 2020), copy format from M3
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 2020

# Data rows 5-14: copy format from column M, then set column N value
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 2198.7

$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 132.7

$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N7").Value = 242.9

$ws.Range("M8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = 203.3

$ws.Range("M9").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N9").Value = 202.8

$ws.Range("M10").Copy()
$ws.Range("N10").PasteSpecial(-4122)
$ws.Range("N10").Value = 284.7

$ws.Range("M11").Copy()
$ws.Range("N11").PasteSpecial(-4122)
$ws.Range("N11").Value = 294.9

$ws.Range("M12").Copy()
$ws.Range("N12").PasteSpecial(-4122)
$ws.Range("N12").Value = 802.5

$ws.Range("M13").Copy()
$ws.Range("N13").PasteSpecial(-4122)
$ws.Range("N13").Value = 28.1

$ws.Range("M14").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Range("N14").Value = 6.8

# Row 15: totals/spacer row, format only, no value
$ws.Range("M15").Copy()
$ws.Range("N15").PasteSpecial(-4122)

# Rows 16-25: percentage rows
$ws.Range("M16").Copy()
$ws.Range("N16").PasteSpecial(-4122)
$ws.Range("N16").Value = 27.4

$ws.Range("M17").Copy()
$ws.Range("N17").PasteSpecial(-4122)
$ws.Range("N17").Value = 17.5

$ws.Range("M18").Copy()
$ws.Range("N18").PasteSpecial(-4122)
$ws.Range("N18").Value = 24.7

$ws.Range("M19").Copy()
$ws.Range("N19").PasteSpecial(-4122)
$ws.Range("N19").Value = 31.5

$ws.Range("M20").Copy()
$ws.Range("N20").PasteSpecial(-4122)
$ws.Range("N20").Value = 30.4

$ws.Range("M21").Copy()
$ws.Range("N21").PasteSpecial(-4122)
$ws.Range("N21").Value = 24.8

$ws.Range("M22").Copy()
$ws.Range("N22").PasteSpecial(-4122)
$ws.Range("N22").Value = 30.7

$ws.Range("M23").Copy()
$ws.Range("N23").PasteSpecial(-4122)
$ws.Range("N23").Value = 30.1

$ws.Range("M24").Copy()
$ws.Range("N24").PasteSpecial(-4122)
$ws.Range("N24").Value = 21.2

$ws.Range("M25").Copy()
$ws.Range("N25").PasteSpecial(-4122)
$ws.Range("N25").Value = 11.6

# Update the active selection to M25, as reflected in the sheetView
$null = $ws.Range("M25").Select()
